# Auto-generated update of leve profit/price figures across multiple sheets.
# Mirrors a scheduled-runner refresh of market board price data.
$wb = $excel.ActiveWorkbook

# ALC (sheet 1), row 17
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 1011.9828
$ws.Range("J17").Value = 835.6491
$ws.Range("L17").Value = 2506.9473
$ws.Range("N17").Value = -2842.9473

# ALC (sheet 1), row 33
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 92
$ws.Range("I33").Value = 75.5
$ws.Range("K33").Value = 75.5
$ws.Range("M33").Value = 153.5

# ALC (sheet 1), row 53
$ws = $wb.Worksheets.Item(1)
$ws.Range("H53").Value = 994.25
$ws.Range("I53").Value = 1280.125
$ws.Range("J53").Value = 422.5
$ws.Range("K53").Value = 1280.125
$ws.Range("L53").Value = 422.5
$ws.Range("M53").Value = -643.125
$ws.Range("N53").Value = -1696.5

# ALC (sheet 1), row 58
$ws = $wb.Worksheets.Item(1)
$ws.Range("H58").Value = 1051.1333
$ws.Range("I58").Value = 315.1
$ws.Range("J58").Value = 2523.2
$ws.Range("K58").Value = 945.3000000000001
$ws.Range("L58").Value = 7569.599999999999
$ws.Range("M58").Value = -795.3000000000001
$ws.Range("N58").Value = -7869.599999999999

# ALC (sheet 1), row 129
$ws = $wb.Worksheets.Item(1)
$ws.Range("H129").Value = 912.0833
$ws.Range("J129").Value = 890.3134
$ws.Range("L129").Value = 2670.9402
$ws.Range("N129").Value = -12670.9402

# ALC (sheet 1), row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 1104.3125
$ws.Range("J132").Value = 866.3333
$ws.Range("L132").Value = 2598.9999
$ws.Range("N132").Value = -7658.9999

# ALC (sheet 1), row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 2557.027
$ws.Range("I138").Value = 2870.5908
$ws.Range("J138").Value = 2097.1333
$ws.Range("K138").Value = 8611.7724
$ws.Range("L138").Value = 6291.3999
$ws.Range("M138").Value = -3471.7724
$ws.Range("N138").Value = -16571.3999

# ARM (sheet 2), row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 4723.298
$ws.Range("I32").Value = 3156.0645
$ws.Range("K32").Value = 3156.0645
$ws.Range("M32").Value = -2869.0645

# ARM (sheet 2), row 74
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 921.375
$ws.Range("I74").Value = 759.6667
$ws.Range("K74").Value = 759.6667
$ws.Range("M74").Value = 114.3333

# ARM (sheet 2), row 77
$ws = $wb.Worksheets.Item(2)
$ws.Range("H77").Value = 921.375
$ws.Range("I77").Value = 759.6667
$ws.Range("K77").Value = 3798.3335
$ws.Range("M77").Value = 569.6665000000003

# ARM (sheet 2), row 122
$ws = $wb.Worksheets.Item(2)
$ws.Range("H122").Value = 1719.6
$ws.Range("I122").Value = 1482
$ws.Range("K122").Value = 4446
$ws.Range("M122").Value = -1996

# BSM (sheet 3), row 100
$ws = $wb.Worksheets.Item(3)
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164

# BSM (sheet 3), row 105
$ws = $wb.Worksheets.Item(3)
$ws.Range("H105").Value = 1854.238
$ws.Range("I105").Value = 1982.2941
$ws.Range("K105").Value = 1982.2941
$ws.Range("M105").Value = -235.2941000000001

# CRP (sheet 4), row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 3908.5557
$ws.Range("I31").Value = 1529.5
$ws.Range("K31").Value = 1529.5
$ws.Range("M31").Value = -1234.5

# CRP (sheet 4), row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 3908.5557
$ws.Range("I34").Value = 1529.5
$ws.Range("K34").Value = 1529.5
$ws.Range("M34").Value = -1327.5

# CRP (sheet 4), row 105
$ws = $wb.Worksheets.Item(4)
$ws.Range("H105").Value = 1502.8334
$ws.Range("I105").Value = 703.5
$ws.Range("K105").Value = 703.5
$ws.Range("M105").Value = 1043.5

# CUL (sheet 5), row 2
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 600
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -826

# CUL (sheet 5), row 38
$ws = $wb.Worksheets.Item(5)
$ws.Range("H38").Value = 286.03705
$ws.Range("I38").Value = 51.75
$ws.Range("J38").Value = 384.6842
$ws.Range("K38").Value = 155.25
$ws.Range("L38").Value = 1154.0526
$ws.Range("M38").Value = 191.75
$ws.Range("N38").Value = -1848.0526

# CUL (sheet 5), row 68
$ws = $wb.Worksheets.Item(5)
$ws.Range("H68").Value = 800
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# CUL (sheet 5), row 71
$ws = $wb.Worksheets.Item(5)
$ws.Range("H71").Value = 800
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# CUL (sheet 5), row 107
$ws = $wb.Worksheets.Item(5)
$ws.Range("H107").Value = 733.1177
$ws.Range("I107").Value = 303
$ws.Range("K107").Value = 909
$ws.Range("M107").Value = 1011

# CUL (sheet 5), row 110
$ws = $wb.Worksheets.Item(5)
$ws.Range("H110").Value = 4999
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 4999
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 14997
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -23177

# CUL (sheet 5), row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 776.86
$ws.Range("J131").Value = 813.4783
$ws.Range("L131").Value = 2440.4349
$ws.Range("N131").Value = -12520.4349

# CUL (sheet 5), row 134
$ws = $wb.Worksheets.Item(5)
$ws.Range("H134").Value = 1394.8846
$ws.Range("I134").Value = 1178.55
$ws.Range("J134").Value = 2116
$ws.Range("K134").Value = 3535.65
$ws.Range("L134").Value = 6348
$ws.Range("M134").Value = 1534.35
$ws.Range("N134").Value = -16488

# GSM (sheet 6), row 122
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# LTW (sheet 7), row 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 3330.6191
$ws.Range("I7").Value = 1752.5
$ws.Range("J7").Value = 8380.6
$ws.Range("K7").Value = 1752.5
$ws.Range("L7").Value = 8380.6
$ws.Range("M7").Value = -1640.5
$ws.Range("N7").Value = -8604.6

# LTW (sheet 7), row 40
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 13333
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 13333
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 13333
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -13605

# LTW (sheet 7), row 55
$ws = $wb.Worksheets.Item(7)
$ws.Range("H55").Value = 526.125
$ws.Range("I55").Value = 489.75
$ws.Range("J55").Value = 562.5
$ws.Range("K55").Value = 489.75
$ws.Range("L55").Value = 562.5
$ws.Range("M55").Value = -316.75
$ws.Range("N55").Value = -908.5

# LTW (sheet 7), row 63
$ws = $wb.Worksheets.Item(7)
$ws.Range("H63").Value = 29110
$ws.Range("J63").Value = 29110
$ws.Range("L63").Value = 29110
$ws.Range("N63").Value = -30608

# LTW (sheet 7), row 66
$ws = $wb.Worksheets.Item(7)
$ws.Range("H66").Value = 29110
$ws.Range("J66").Value = 29110
$ws.Range("L66").Value = 87330
$ws.Range("N66").Value = -94818

# LTW (sheet 7), row 74
$ws = $wb.Worksheets.Item(7)
$ws.Range("H74").Value = 20000000
$ws.Range("J74").Value = 20000000
$ws.Range("L74").Value = 20000000
$ws.Range("N74").Value = -20001996

# LTW (sheet 7), row 77
$ws = $wb.Worksheets.Item(7)
$ws.Range("H77").Value = 20000000
$ws.Range("J77").Value = 20000000
$ws.Range("L77").Value = 60000000
$ws.Range("N77").Value = -60009984

# LTW (sheet 7), row 122
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 3883.6155
$ws.Range("J122").Value = 4299
$ws.Range("L122").Value = 12897
$ws.Range("N122").Value = -17797

# LTW (sheet 7), row 126
$ws = $wb.Worksheets.Item(7)
$ws.Range("H126").Value = 3330.6191
$ws.Range("I126").Value = 1752.5
$ws.Range("J126").Value = 8380.6
$ws.Range("K126").Value = 5257.5
$ws.Range("L126").Value = 25141.8
$ws.Range("M126").Value = -2787.5
$ws.Range("N126").Value = -30081.8

# LTW (sheet 7), row 136
$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 3150.4
$ws.Range("I136").Value = 1723.6875
$ws.Range("K136").Value = 5171.0625
$ws.Range("M136").Value = -2621.0625

# WVR (sheet 8), row 113
$ws = $wb.Worksheets.Item(8)
$ws.Range("H113").Value = 629.35297
$ws.Range("I113").Value = 470.69232
$ws.Range("J113").Value = 1145
$ws.Range("K113").Value = 1412.07696
$ws.Range("L113").Value = 3435
$ws.Range("M113").Value = 757.9230400000001
$ws.Range("N113").Value = -7775

# WVR (sheet 8), row 132
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 5346.2383
$ws.Range("I132").Value = 1174
$ws.Range("K132").Value = 3522
$ws.Range("M132").Value = -992

# WVR (sheet 8), row 136
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 17923568
$ws.Range("I136").Value = 27780072
$ws.Range("J136").Value = 2653.2727
$ws.Range("K136").Value = 83340216
$ws.Range("L136").Value = 7959.8181
$ws.Range("M136").Value = -83337666
$ws.Range("N136").Value = -13059.8181
